# Section 26 - Using Excel's LEN() Function
# Adds a helper "Num Chars" column (H) that computes LEN(A) for each SKU,
# and updates the RIGHT() formula in column G to branch on that length so
# it still works for SKUs that aren't exactly 8 characters long. Also
# back-fills the LEFT/MID/RIGHT helper formulas for rows 18-26, which
# previously had no formulas at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LEFT RIGHT MID Functions")
$wb.Activate()
$ws.Activate()

# --- Row 3: plain (non-shared) formulas, matching the existing E3/F3 style ---
$ws.Range("H3").Formula = "=LEN(A3)"
$ws.Range("G3").Formula = "=IF(H3 = 8, RIGHT(A3,2), RIGHT(A3,4))"

# --- Row 4: becomes the shared-formula "master" cell for G4:G26 / H4:H26 ---
$ws.Range("H4").Formula = "=LEN(A4)"
$ws.Range("G4").Formula = "=IF(H4 = 8, RIGHT(A4,2), RIGHT(A4,4))"

# --- Fill G/H down through the rest of the table (rows 5-26) ---
for ($r = 5; $r -le 26; $r++) {
    $ws.Range("H$r").Formula = "=LEN(A$r)"
    $ws.Range("G$r").Formula = "=IF(H$r = 8, RIGHT(A$r,2), RIGHT(A$r,4))"
}

# --- Rows 18-26 previously had empty E/F/G cells; back-fill E (LEFT) and
#     F (MID) the same way the rest of the table already works ---
for ($r = 18; $r -le 26; $r++) {
    $ws.Range("E$r").Formula = "=LEFT(A$r,3)"
    $ws.Range("F$r").Formula = "=MID(A$r,4,3)"
}

# --- Selection moved to J20 as part of the edit ---
$ws.Range("J20").Select()

$excel.CalculateFullRebuild()
